$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("E2").Value = 3
$ws.Range("F2").Value = 1
$ws.Range("G2").Value = 165.404724
$ws.Range("H2").Value = 496.214172
$ws.Range("I2").Value = 0.3557141051771751
$ws.Range("J2").Value = 0.355714105177175
$ws.Range("M2").Value = 174.1282373333333
$ws.Range("N2").Value = 522.384712
$ws.Range("O2").Value = 0.985625830323027
$ws.Range("P2").Value = 0.985625830323027
$ws.Range("Q2").Value = 28801.6330367265
$ws.Range("R2").Value = 259214.6973305385
$ws.Range("S2").Value = 0.3506010102728657
$ws.Range("T2").Value = 0.3506010102728657
$ws.Range("E3").Value = 3
$ws.Range("F3").Value = 1
$ws.Range("G3").Value = 165.404724
$ws.Range("H3").Value = 496.214172
$ws.Range("I3").Value = 0.3557141051771751
$ws.Range("J3").Value = 0.355714105177175
$ws.Range("O3").Value = 0.003686901313133159
$ws.Range("P3").Value = 0.003686901313133159
$ws.Range("Q3").Value = 107.737414540652
$ws.Range("R3").Value = 969.6367308658681
$ws.Range("S3").Value = 0.001311482801477713
$ws.Range("T3").Value = 0.001311482801477713
$ws.Range("E4").Value = 3
$ws.Range("F4").Value = 1
$ws.Range("G4").Value = 165.404724
$ws.Range("H4").Value = 496.214172
$ws.Range("I4").Value = 0.3557141051771751
$ws.Range("J4").Value = 0.355714105177175
$ws.Range("M4").Value = 1.888095
$ws.Range("N4").Value = 5.664285
$ws.Range("O4").Value = 0.01068726836383999
$ws.Range("P4").Value = 0.01068726836383999
$ws.Range("Q4").Value = 312.29983236078
$ws.Range("R4").Value = 2810.69849124702
$ws.Range("S4").Value = 0.003801612102831675
$ws.Range("T4").Value = 0.003801612102831675
$ws.Range("I5").Value = 0.2830704861820888
$ws.Range("J5").Value = 0.2830704861820888
$ws.Range("M5").Value = 174.1282373333333
$ws.Range("N5").Value = 522.384712
$ws.Range("O5").Value = 0.985625830323027
$ws.Range("P5").Value = 0.985625830323027
$ws.Range("Q5").Value = 22919.78908872187
$ws.Range("R5").Value = 206278.1017984968
$ws.Range("S5").Value = 0.2790015829831642
$ws.Range("T5").Value = 0.2790015829831642
$ws.Range("I6").Value = 0.2830704861820888
$ws.Range("J6").Value = 0.2830704861820888
$ws.Range("O6").Value = 0.003686901313133159
$ws.Range("P6").Value = 0.003686901313133159
$ws.Range("S6").Value = 0.001043652947213985
$ws.Range("T6").Value = 0.001043652947213985
$ws.Range("I7").Value = 0.2830704861820888
$ws.Range("J7").Value = 0.2830704861820888
$ws.Range("M7").Value = 1.888095
$ws.Range("N7").Value = 5.664285
$ws.Range("O7").Value = 0.01068726836383999
$ws.Range("P7").Value = 0.01068726836383999
$ws.Range("Q7").Value = 248.52223764621
$ws.Range("R7").Value = 2236.70013881589
$ws.Range("S7").Value = 0.003025250251710644
$ws.Range("T7").Value = 0.003025250251710644
$ws.Range("G8").Value = 167.962794
$ws.Range("H8").Value = 503.888382
$ws.Range("I8").Value = 0.3612154086407362
$ws.Range("J8").Value = 0.3612154086407361
$ws.Range("M8").Value = 174.1282373333333
$ws.Range("N8").Value = 522.384712
$ws.Range("O8").Value = 0.985625830323027
$ws.Range("P8").Value = 0.985625830323027
$ws.Range("Q8").Value = 29247.06525680178
$ws.Range("R8").Value = 263223.587311216
$ws.Range("S8").Value = 0.3560232370669971
$ws.Range("T8").Value = 0.3560232370669971
$ws.Range("G9").Value = 167.962794
$ws.Range("H9").Value = 503.888382
$ws.Range("I9").Value = 0.3612154086407362
$ws.Range("J9").Value = 0.3612154086407361
$ws.Range("O9").Value = 0.003686901313133159
$ws.Range("P9").Value = 0.003686901313133159
$ws.Range("Q9").Value = 109.403629636262
$ws.Range("R9").Value = 984.6326667263581
$ws.Range("S9").Value = 0.001331765564441461
$ws.Range("T9").Value = 0.001331765564441461
$ws.Range("G10").Value = 167.962794
$ws.Range("H10").Value = 503.888382
$ws.Range("I10").Value = 0.3612154086407362
$ws.Range("J10").Value = 0.3612154086407361
$ws.Range("M10").Value = 1.888095
$ws.Range("N10").Value = 5.664285
$ws.Range("O10").Value = 0.01068726836383999
$ws.Range("P10").Value = 0.01068726836383999
$ws.Range("Q10").Value = 317.12971153743
$ws.Range("R10").Value = 2854.167403836871
$ws.Range("S10").Value = 0.003860406009297676
$ws.Range("T10").Value = 0.003860406009297675